# Sar1 and Arf1 AH NEES data added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 46-47 append two more AH sequences (Sar1, Arf1) below the
# existing AA_seq / Protein_Name / AH# table.
# Cell values are written in this specific order so that the shared-string
# table ends up populated as: Sar1, Arf1, MGNIFANLFKGLFGKKE, MAGWDIFGWFRDVLASLGLWNKH
$ws.Range("B46").Value = "Sar1"
$ws.Range("B47").Value = "Arf1"
$ws.Range("A47").Value = "MGNIFANLFKGLFGKKE"
$ws.Range("A46").Value = "MAGWDIFGWFRDVLASLGLWNKH"

$ws.Range("C46").Value = 1
$ws.Range("C47").Value = 1

# Reflect the scrolled viewport / new active selection from the edit session
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("B46").Select()
